$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / country-name updates (shared-string table reorder equivalent) ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 26 de Junio de 2020 a las 01:44"
$ws.Cells.Item(31, 1).Value = "Argentina"
$ws.Cells.Item(32, 1).Value = "Indonesia"
$ws.Cells.Item(33, 1).Value = "Paises Bajos"
$ws.Cells.Item(51, 1).Value = "Nigeria"
$ws.Cells.Item(52, 1).Value = "Armenia"
$ws.Cells.Item(53, 1).Value = "Israel"
$ws.Cells.Item(166, 1).Value = "Guyana"
$ws.Cells.Item(167, 1).Value = "Angola"
$ws.Cells.Item(202, 1).Value = "Fiyi"
$ws.Cells.Item(203, 1).Value = "Dominica"
$ws.Cells.Item(208, 1).Value = "Groenlandia"
$ws.Cells.Item(209, 1).Value = "Islas Malvinas"
$ws.Cells.Item(211, 1).Value = "Seychelles"
$ws.Cells.Item(212, 1).Value = "Montserrat"

# --- Numeric statistic updates ---
$ws.Cells.Item(4, 2).Value = 2499127
$ws.Cells.Item(4, 3).Value = 36573
$ws.Cells.Item(4, 4).Value = 1047894
$ws.Cells.Item(4, 5).Value = 1326368
$ws.Cells.Item(4, 7).Value = 584
$ws.Cells.Item(4, 8).Value = 124865
$ws.Cells.Item(5, 2).Value = 1233147
$ws.Cells.Item(5, 3).Value = 40673
$ws.Cells.Item(5, 5).Value = 528185
$ws.Cells.Item(5, 7).Value = 1180
$ws.Cells.Item(5, 8).Value = 55054
$ws.Cells.Item(15, 2).Value = 193785
$ws.Cells.Item(15, 3).Value = 531
$ws.Cells.Item(15, 5).Value = 7973
$ws.Cells.Item(22, 2).Value = 102622
$ws.Cells.Item(22, 3).Value = 380
$ws.Cells.Item(22, 4).Value = 65425
$ws.Cells.Item(22, 5).Value = 28693
$ws.Cells.Item(22, 7).Value = 20
$ws.Cells.Item(22, 8).Value = 8504
$ws.Cells.Item(31, 2).Value = 52457
$ws.Cells.Item(31, 3).Value = 2606
$ws.Cells.Item(31, 4).Value = 14788
$ws.Cells.Item(31, 5).Value = 36519
$ws.Cells.Item(31, 7).Value = 34
$ws.Cells.Item(31, 8).Value = 1150
$ws.Cells.Item(32, 2).Value = 50187
$ws.Cells.Item(32, 3).Value = 1178
$ws.Cells.Item(32, 4).Value = 20449
$ws.Cells.Item(32, 5).Value = 27118
$ws.Cells.Item(32, 7).Value = 47
$ws.Cells.Item(32, 8).Value = 2620
$ws.Cells.Item(33, 2).Value = 49914
$ws.Cells.Item(33, 3).Value = 110
$ws.Cells.Item(33, 4).Value = 0
$ws.Cells.Item(33, 5).Value = 0
$ws.Cells.Item(33, 7).Value = 3
$ws.Cells.Item(33, 8).Value = 6100
$ws.Cells.Item(51, 2).Value = 22614
$ws.Cells.Item(51, 3).Value = 594
$ws.Cells.Item(51, 4).Value = 7822
$ws.Cells.Item(51, 5).Value = 14243
$ws.Cells.Item(51, 7).Value = 7
$ws.Cells.Item(51, 8).Value = 549
$ws.Cells.Item(52, 2).Value = 22488
$ws.Cells.Item(52, 3).Value = 771
$ws.Cells.Item(52, 4).Value = 11335
$ws.Cells.Item(52, 5).Value = 10756
$ws.Cells.Item(52, 7).Value = 11
$ws.Cells.Item(52, 8).Value = 397
$ws.Cells.Item(53, 2).Value = 22400
$ws.Cells.Item(53, 3).Value = 356
$ws.Cells.Item(53, 4).Value = 16007
$ws.Cells.Item(53, 5).Value = 6084
$ws.Cells.Item(53, 7).Value = 1
$ws.Cells.Item(53, 8).Value = 309
$ws.Cells.Item(69, 2).Value = 10870
$ws.Cells.Item(69, 3).Value = 93
$ws.Cells.Item(69, 5).Value = 2876
$ws.Cells.Item(70, 2).Value = 8984
$ws.Cells.Item(70, 3).Value = 95
$ws.Cells.Item(70, 4).Value = 3806
$ws.Cells.Item(70, 5).Value = 4622
$ws.Cells.Item(70, 7).Value = 8
$ws.Cells.Item(70, 8).Value = 556
$ws.Cells.Item(73, 2).Value = 8334
$ws.Cells.Item(73, 3).Value = 170
$ws.Cells.Item(73, 4).Value = 3487
$ws.Cells.Item(73, 5).Value = 4787
$ws.Cells.Item(73, 7).Value = 2
$ws.Cells.Item(73, 8).Value = 60
$ws.Cells.Item(86, 2).Value = 5087
$ws.Cells.Item(86, 3).Value = 131
$ws.Cells.Item(86, 4).Value = 2270
$ws.Cells.Item(86, 5).Value = 2777
$ws.Cells.Item(86, 7).Value = 1
$ws.Cells.Item(86, 8).Value = 40
$ws.Cells.Item(94, 2).Value = 3739
$ws.Cells.Item(94, 3).Value = 220
$ws.Cells.Item(94, 4).Value = 1225
$ws.Cells.Item(94, 5).Value = 2395
$ws.Cells.Item(94, 7).Value = 3
$ws.Cells.Item(94, 8).Value = 119
$ws.Cells.Item(137, 2).Value = 907
$ws.Cells.Item(137, 3).Value = 5
$ws.Cells.Item(137, 4).Value = 818
$ws.Cells.Item(137, 5).Value = 63
$ws.Cells.Item(144, 2).Value = 711
$ws.Cells.Item(144, 3).Value = 1
$ws.Cells.Item(144, 4).Value = 214
$ws.Cells.Item(144, 5).Value = 484
$ws.Cells.Item(151, 2).Value = 588
$ws.Cells.Item(151, 3).Value = 5
$ws.Cells.Item(151, 4).Value = 394
$ws.Cells.Item(151, 5).Value = 180
$ws.Cells.Item(152, 2).Value = 551
$ws.Cells.Item(152, 3).Value = 21
$ws.Cells.Item(152, 4).Value = 128
$ws.Cells.Item(152, 5).Value = 417
$ws.Cells.Item(157, 2).Value = 373
$ws.Cells.Item(157, 3).Value = 16
$ws.Cells.Item(157, 4).Value = 176
$ws.Cells.Item(157, 5).Value = 187
$ws.Cells.Item(166, 2).Value = 215
$ws.Cells.Item(166, 3).Value = 6
$ws.Cells.Item(166, 4).Value = 108
$ws.Cells.Item(166, 5).Value = 95
$ws.Cells.Item(166, 8).Value = 12
$ws.Cells.Item(167, 2).Value = 212
$ws.Cells.Item(167, 3).Value = 15
$ws.Cells.Item(167, 4).Value = 81
$ws.Cells.Item(167, 5).Value = 121
$ws.Cells.Item(167, 8).Value = 10
$ws.Cells.Item(194, 4).Value = 29
$ws.Cells.Item(194, 5).Value = 0
$ws.Cells.Item(211, 4).Value = 11
$ws.Cells.Item(211, 8).Value = 0
$ws.Cells.Item(212, 4).Value = 10
$ws.Cells.Item(212, 8).Value = 1
